# Estructura de actividades RGG 20240509
# Applies the activity-log updates to "actividades" (new rows for 2024-05-09)
# and documents the "control_digit3" work item on the
# "descriptivo- detalle a migrar f" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "actividades": insert 4 new rows at the top of the log
# (2024-05-09 activities), pushing the existing rows down.
# ---------------------------------------------------------------
$act = $wb.Worksheets.Item("actividades")

$act.Rows("2:5").Insert()

# Copy the date-column format (style) from the row that used to be
# row 2 (now row 6) down onto the freshly inserted date cells.
$act.Range("A6").Copy()
$act.Range("A2:A5").PasteSpecial(-4122)

# Copy the description-column format (vertical-centered style, like
# the one used further down in the sheet) onto the new D column cells.
$act.Range("D25").Copy()
$act.Range("D2:D5").PasteSpecial(-4122)

$act.Range("A2:A5").Value = 45421

$act.Range("B2").Value = 1
$act.Range("B3").Value = 2
$act.Range("B4").Value = 3
$act.Range("B5").Value = 4

$act.Range("C2").Value = "SPOOLER"
$act.Range("C3").Value = "SPOOLER"
$act.Range("C4").Value = "SPOOLER"
$act.Range("C5").Value = "DOUMENTACION"

$act.Range("D2").Value = " procedimiento almacenado se genero la aplicación y se versiono en la ruta de doucmentaciones carpeta (C:\pc\raul\Net\migracion_spooler\00-Documentacion\db\bdoracle)"
$act.Range("D3").Value = "Aplicación se modificó la aplicación se unieron los los 3 sql  , se modificó la función principal de sql, agregan para metro opcional. Para el llenado del 3er sql."
$act.Range("D4").Value = "Se investigo la ejecucion función los cuales , se ajustaron envia error al llenar la tabla “System.Exception: 'Size must be set.'”"

$act.Range("D5").Value = "Se continuo con al documentación de control_digit3"
$prefixLen = ("Se continuo con al documentación de ").Length
$fullLen = ("Se continuo con al documentación de control_digit3").Length
$runChars = $act.Range("D5").Characters($prefixLen + 1, $fullLen - $prefixLen)
$runChars.Font.Name = "Calibri"
$runChars.Font.Size = 11
$runChars.Font.Color = 0

# ---------------------------------------------------------------
# Sheet "descriptivo- detalle a migrar f": add the control_digit3
# function block (rows 87-93) after the existing control_digit
# block that ends at row 86.
# ---------------------------------------------------------------
$det = $wb.Worksheets.Item("descriptivo- detalle a migrar f")

# Row 87 (section header row) - copy the style skeleton from row 80
# (A:H only, column I intentionally left untouched/absent).
$det.Range("A80:H80").Copy()
$det.Range("A87:H87").PasteSpecial(-4122)
$det.Rows("87").RowHeight = 30

$det.Range("A87").Value = 188
$det.Range("C87").Value = "control_digit3"
$det.Range("D87").Value = "General un excel - con inf. BD"
$det.Range("E87").Value = "correo,ftp,xml, excel"
$det.Range("F87").Value = "control_digit3"
$det.Range("G87").Value = 12
$det.Range("H87").Value = 4

# Row 88 - detail row, same shape as row 81 (D,E,F,I).
$det.Range("D81:I81").Copy()
$det.Range("D88:I88").PasteSpecial(-4122)
$det.Range("F88").Value = "init_var"
$det.Range("I88").Value = 1

# Row 89 - detail row, same shape as row 83 (D,E,F,H,I,K no style).
$det.Range("D83:K83").Copy()
$det.Range("D89:K89").PasteSpecial(-4122)
$det.Range("F89").Value = "ftp_sucursal_cargar"
$det.Range("H89").Value = 1
$det.Range("I89").Value = 1.1
$det.Range("K89").Value = "valida las sucursales que contenga dicho proceso  de pertenece a unidad  ftp_sucursal.bas"

# Row 90 - detail row, same shape as row 86 (D,E,F,H,I, no K).
$det.Range("D86:I86").Copy()
$det.Range("D90:I90").PasteSpecial(-4122)
$det.Range("F90").Value = "SQL_INSERT"
$det.Range("H90").Value = 4
$det.Range("I90").Value = 2

# Row 91 - detail row, same shape as row 82 (D,E,F,H,I,K styled).
$det.Range("D82:K82").Copy()
$det.Range("D91:K91").PasteSpecial(-4122)
$det.Range("F91").Value = "log_SQL"
$det.Range("H91").Value = 1
$det.Range("I91").Value = 3
$det.Range("K91").Value = "log- bitacora"

# Row 92 - detail row, same shape as row 83 (D,E,F,H,I,K no style).
$det.Range("D83:K83").Copy()
$det.Range("D92:K92").PasteSpecial(-4122)
$det.Range("F92").Value = "validar_evidencia"
$det.Range("H92").Value = 0
$det.Range("I92").Value = 4
$det.Range("K92").Value = "NA"

# Row 93 - detail row, same shape as row 86 (D,E,F,H,I, no K).
$det.Range("D86:I86").Copy()
$det.Range("D93:I93").PasteSpecial(-4122)
$det.Range("F93").Value = "SQL_DIGIT"
$det.Range("H93").Value = 2
$det.Range("I93").Value = 5

# Update the "correo" -> "correo,excel" notes on the two other
# section header rows (the notification channels used grew to
# include excel export).
$det.Range("E71").Value = "correo,excel"
$det.Range("E72").ClearContents()
$det.Range("E80").Value = "correo,excel"

# ---------------------------------------------------------------
# Window / selection bookkeeping to match the saved state: the
# "actividades" sheet becomes the active tab, selection on A6, and
# the "descriptivo..." sheet keeps its frozen-pane view with the
# new bottom row selected.
# ---------------------------------------------------------------
$det.Range("G87").Select()
$act.Activate()
$act.Range("A6").Select()
